# refactor add company automation script
# - Row 3 (Login.login test data): replace phone/verification-code numbers
#   with the new placeholder creds "james" / "123456".
# - Row 6 (addCompany test data): replace the stray "forTesting" value with
#   the real "company details " payload.
# - Move the active selection to E6 (last edited cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "james"
$ws.Range("E3").Value = "123456"
$ws.Range("E6").Value = "company details "

[void]$ws.Range("E6").Select()
